$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---
$ws1.Range("D2").Value = 13
$ws1.Range("H2").Value = 7.92
$ws1.Range("L2").Value = 0.91

$ws1.Range("H3").Value = 5.77
$ws1.Range("L3").Value = 1.15

$ws1.Range("H4").Value = 6.36
$ws1.Range("L4").Value = 0.92

$ws1.Range("H5").Value = 4.82
$ws1.Range("L5").Value = 1.06

$ws1.Range("H6").Value = 3.82
$ws1.Range("L6").Value = 1.07

$ws1.Range("D7").Value = 12
$ws1.Range("H7").Value = 3.06
$ws1.Range("L7").Value = 1.16

$ws1.Range("D8").Value = 11
$ws1.Range("H8").Value = 2.25

$ws1.Range("D9").Value = 12
$ws1.Range("H9").Value = 1.14
$ws1.Range("J9").Value = "Normal"
$ws1.Range("L9").Value = 1.03

$ws1.Range("H10").Value = 0.15
$ws1.Range("L10").Value = 1.13

$ws1.Range("L11").Value = 1.16

$ws1.Range("L12").Value = 0.84

$ws1.Range("L13").Value = 0.95

$ws1.Range("L14").Value = 0.82

$ws1.Range("L15").Value = 1.11

$ws1.Range("L16").Value = 0.96

$ws1.Range("L17").Value = 1.2

# --- Summary sheet ---
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "192"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "101"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "53"
